# Applies the "Saldo" export refresh: new accounts added, several balances
# updated (causing re-sort by descending Saldo), and a few stale rows
# removed. Operations are performed bottom-to-top so earlier row numbers
# stay valid as the sheet reflows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 134: SILVIA (005547703) old low balance 75.37 -> removed
#     (she now appears near the top with an updated balance, see below) ---
$ws.Rows.Item(134).Delete() | Out-Null

# --- Row 24: RAPHAELA (005366255) old balance 844.05 -> removed
#     (replaced by an updated row inserted just above, see below) ---
$ws.Rows.Item(24).Delete() | Out-Null

# --- Insert RAPHAELA (005366255) with updated balance 868.76 before row 23 ---
$ws.Rows.Item(23).Insert() | Out-Null
$ws.Cells.Item(23, 1).NumberFormat = "@"
$ws.Cells.Item(23, 1).Value = "005366255"
$ws.Cells.Item(23, 2).Value = "RAPHAELA"
$ws.Cells.Item(23, 3).Value = 868.76

# --- Rows 17-18: GUSTAVO (004565108) / TATYANA (004466342) -> removed ---
$ws.Rows.Item(17).Resize(2).Delete() | Out-Null

# --- Rows 14-15: BRUNO (004515341) / GABRIELA (004431546) -> replaced by a
#     single updated row for TATIANA (005366671) ---
$ws.Rows.Item(15).Delete() | Out-Null
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Value = "005366671"
$ws.Cells.Item(14, 2).Value = "TATIANA"
$ws.Cells.Item(14, 3).Value = 6600

# --- Rows 9-11: SISSI (004315417) / HENRIQUE (004479463) / HENRIQUE
#     (004497875) -> removed ---
$ws.Rows.Item(9).Resize(3).Delete() | Out-Null

# --- Rows 4-7: PEDRO (005081833) / MSD (004526450) / ADSON (004404342) /
#     NATALIA (004482102) -> replaced by PEDRO's updated balance plus a new
#     SILVIA (005547703) row with her updated balance ---
$ws.Rows.Item(6).Resize(2).Delete() | Out-Null
$ws.Cells.Item(4, 3).Value = 74795.79
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "005547703"
$ws.Cells.Item(5, 2).Value = "SILVIA"
$ws.Cells.Item(5, 3).Value = 59315.58

# --- Insert two brand-new top accounts: DIMITRI (005529100) and
#     FELIPE (005135532) ---
$ws.Rows.Item(2).Resize(2).Insert() | Out-Null
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "005529100"
$ws.Cells.Item(2, 2).Value = "DIMITRI"
$ws.Cells.Item(2, 3).Value = 242497.72
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "005135532"
$ws.Cells.Item(3, 2).Value = "FELIPE"
$ws.Cells.Item(3, 3).Value = 215000
